{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list to use impact-focused\n// accomplishment statements (per commit: \"Fix Key Achievements to use proper\n// accomplishment statements\"). The same bullet wording also appears earlier in\n// the PROFESSIONAL EXPERIENCE section, so we must only touch the occurrence\n// that lives under the \"KEY ACHIEVEMENTS AND IMPACT\" heading.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading paragraph.\nlet headingIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error('Could not find \"KEY ACHIEVEMENTS AND IMPACT\" heading');\n}\n\n// The exact (old) bullet text we expect to find, in order, after the heading.\nconst oldBullets = [\n  \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n  \"\\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n  \"\\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n  \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n  \"\\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n  \"\\u2022 Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\",\n];\n\n// Find these bullets as consecutive paragraphs starting right after the heading\n// (allowing for an intervening \"Impact\" sub-heading paragraph).\nlet start = -1;\nfor (let i = headingIndex + 1; i < items.length - oldBullets.length + 1; i++) {\n  let matches = true;\n  for (let j = 0; j < oldBullets.length; j++) {\n    if (items[i + j].text !== oldBullets[j]) {\n      matches = false;\n      break;\n    }\n  }\n  if (matches) {\n    start = i;\n    break;\n  }\n}\nif (start === -1) {\n  throw new Error(\"Could not locate the expected Key Achievements bullet block\");\n}\n\n// New bullet text: 4 rewritten bullets, replacing the old 6. The final\n// paragraph (old bullet index 5, \"Provided expert testimony...\") is removed\n// entirely; its slot is reused to host the rewritten \"Platform impact\" bullet\n// while the true 6th paragraph is deleted.\nconst newBullets = [\n  \"\\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  \"\\u2022 178% accuracy improvement in racial classification algorithms\",\n  \"\\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  \"\\u2022 $4.7M savings enabled nonprofit access\",\n  \"\\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\",\n];\n\nfor (let j = 0; j < newBullets.length; j++) {\n  items[start + j].insertText(newBullets[j], Word.InsertLocation.replace);\n}\n\n// Delete the now-surplus 6th paragraph (old \"Provided expert testimony...\" slot).\nitems[start + newBullets.length].delete();\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list to use impact-focused\n# accomplishment statements (per commit: \"Fix Key Achievements to use proper\n# accomplishment statements\"). The same bullet wording also appears earlier in\n# the PROFESSIONAL EXPERIENCE section, so we must only touch the occurrence\n# that lives under the \"KEY ACHIEVEMENTS AND IMPACT\" heading.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading paragraph.\n$headingIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading\"\n}\n\n# The exact (old) bullet text we expect to find, in order, after the heading.\n$oldBullets = @(\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\",\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\",\n    \"\u2022 Provided expert testimony and press briefings on electoral data integrity and demographic modeling accuracy\"\n)\n\n# Find these bullets as consecutive paragraphs starting right after the heading\n# (allowing for an intervening \"Impact\" sub-heading paragraph).\n$start = -1\nfor ($i = $headingIndex + 1; $i -le ($count - $oldBullets.Count + 1); $i++) {\n    $matches = $true\n    for ($j = 0; $j -lt $oldBullets.Count; $j++) {\n        $paraText = $d.Paragraphs.Item($i + $j).Range.Text.Trim()\n        if ($paraText -ne $oldBullets[$j]) {\n            $matches = $false\n            break\n        }\n    }\n    if ($matches) {\n        $start = $i\n        break\n    }\n}\nif ($start -eq -1) {\n    throw \"Could not locate the expected Key Achievements bullet block\"\n}\n\n# New bullet text: 4 rewritten bullets, replacing the old 6. The final\n# paragraph (old bullet index 5, \"Provided expert testimony...\") is removed\n# entirely; its slot is reused to host the rewritten \"Platform impact\" bullet\n# while the true 6th paragraph is deleted.\n$newBullets = @(\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\",\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n    \"\u2022 `$4.7M savings enabled nonprofit access\",\n    \"\u2022 Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations\"\n)\n\nfor ($j = 0; $j -lt $newBullets.Count; $j++) {\n    $d.Paragraphs.Item($start + $j).Range.Text = $newBullets[$j]\n}\n\n# Delete the now-surplus 6th paragraph (old \"Provided expert testimony...\" slot).\n$d.Paragraphs.Item($start + $newBullets.Count).Range.Delete()\n"}
